$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10:E14").Value = $true
